{"js": "// Apply the \"FONTE 70A LITE\" listing update: new URL/name/price/store/type/place.\n// Each field lives in its own paragraph as a single \"Label: value\" run, so we\n// find the old full-paragraph text and replace it in place (format-preserving).\nconst replacements = [\n  [\n    \"URL: https://produto.mercadolivre.com.br/MLB-3725911039-fonte-storm-lite-70a-jfa-_JM\",\n    \"URL: https://produto.mercadolivre.com.br/MLB-4834679022-fonte-automotiva-jfa-storm-lite-70a-bivolt-carregador-som-_JM\",\n  ],\n  [\n    \"Nome: Fonte Storm Lite 70a - Jfa\",\n    \"Nome: Fonte Automotiva Jfa Storm Lite 70a Bivolt Carregador Som \",\n  ],\n  [\"Pre\u00e7o: 425.0\", \"Pre\u00e7o: 276.51\"],\n  [\"Pre\u00e7o Previsto: 434.42\", \"Pre\u00e7o Previsto: 408.73\"],\n  [\"Loja: ULTRA COMERCE\", \"Loja: LM.BRASIL\"],\n  [\"Tipo: Premium\", \"Tipo: Cl\u00e1ssico\"],\n  [\"Lugar: Bra\u00e7o do Norte, Santa Catarina.\", \"Lugar: Regente Feij\u00f3, S\u00e3o Paulo.\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"FONTE 70A LITE\" listing update: new URL/name/price/store/type/place.\n# Each field lives in its own paragraph as a single \"Label: value\" run, so we\n# Find & Replace the old full-paragraph text with the new one, one field at a time.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"URL: https://produto.mercadolivre.com.br/MLB-3725911039-fonte-storm-lite-70a-jfa-_JM\",\n      \"URL: https://produto.mercadolivre.com.br/MLB-4834679022-fonte-automotiva-jfa-storm-lite-70a-bivolt-carregador-som-_JM\"),\n    @(\"Nome: Fonte Storm Lite 70a - Jfa\",\n      \"Nome: Fonte Automotiva Jfa Storm Lite 70a Bivolt Carregador Som \"),\n    @(\"Pre\u00e7o: 425.0\", \"Pre\u00e7o: 276.51\"),\n    @(\"Pre\u00e7o Previsto: 434.42\", \"Pre\u00e7o Previsto: 408.73\"),\n    @(\"Loja: ULTRA COMERCE\", \"Loja: LM.BRASIL\"),\n    @(\"Tipo: Premium\", \"Tipo: Cl\u00e1ssico\"),\n    @(\"Lugar: Bra\u00e7o do Norte, Santa Catarina.\", \"Lugar: Regente Feij\u00f3, S\u00e3o Paulo.\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
